$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Mark the title picture's run as NoProofing -> adds <w:noProof/>
#    to that run's rPr (first paragraph / first run of the document).
# ------------------------------------------------------------------
$picPara = $d.Paragraphs.Item(1)
$picPara.Range.NoProofing = $true

# ------------------------------------------------------------------
# 2) "Hours: 17682" -> "Hours: 1" + a new complex-script run "9004"
#    (rFonts hint=cs, rtl) appended right after it, leaving the
#    "Hours: " run untouched.
# ------------------------------------------------------------------
$hoursRng = $d.Content
$hoursRng.Find.Execute("17682", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hoursRng.Text = "1"

$hoursXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
  '<w:r><w:rPr><w:rFonts w:hint="cs"/><w:sz w:val="40"/><w:szCs w:val="40"/><w:rtl/><w:lang w:val="en-US"/></w:rPr><w:t>9004</w:t></w:r>' +
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$hoursRng.InsertXML($hoursXml)

# ------------------------------------------------------------------
# 3) "Cycles: 9585" -> "Cycles: 10125", with the number run losing
#    its <w:lang/> while the "Cycles: " run keeps its own formatting.
# ------------------------------------------------------------------
$cyclesRng = $d.Content
$cyclesRng.Find.Execute("9585", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
# Keep a 1-character placeholder so the range stays anchored right
# after "Cycles: " (a fully-collapsed range re-targets InsertXML to
# the end of the document instead of the intended insertion point).
$cyclesRng.Text = "X"

$cyclesXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
  '<w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>10125</w:t></w:r>' +
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cyclesRng.InsertXML($cyclesXml)

# Drop the "X" placeholder, restoring "Cycles: " as the sole prefix run.
$cleanupRng = $d.Content
$cleanupRng.Find.Execute("Cycles: X", $true, $false, $false, $false, $false, $true, 1, $false, "Cycles: ", 2) | Out-Null
